# Resolved Issue #8 Duplicate lines in BOM
#
# The BOM listed some identical parts (same description / mfg / vendor /
# cost) on several separate rows, one per designator or small designator
# group, instead of a single row with a combined designator list and a
# summed quantity. This script consolidates those duplicate rows:
#
#   - C50 (qty 1), C2 (qty 1), C7/C8/C9/C10 (qty 4)  -> one row, qty 6,
#     designator "C50, C2, C7, C8, C9, C10"
#   - R3 (qty 1), R8/R9 (qty 2)                      -> one row, qty 3,
#     designator "R3, R8, R9"
#
# It also removes a stray placeholder row (all "*" cells, designator
# "R10, R11") that was left over as junk data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Consolidate the C1608X7R1E104K080AA (.1uF) capacitor rows ---------
# Row 8 = C50 (qty 1), row 9 = C2 (qty 1), row 10 = C7,C8,C9,C10 (qty 4).
# Keep row 8, combine qty + designators, then delete rows 9 and 10.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 9).Value = "C50, C2, C7, C8, C9, C10"
$ws.Rows("9:10").Delete()

# --- Consolidate the ERJ-3EKF5601V (5.6 kohm) resistor rows -------------
# After the deletion above, the old row 18 (R3, qty 1) is now row 16 and
# the old row 19 (R8, R9, qty 2) is now row 17.
# Keep row 16, combine qty + designators, then delete row 17.
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 9).Value = "R3, R8, R9"
$ws.Rows("17:17").Delete()

# --- Remove the leftover junk placeholder row ----------------------------
# After the two deletions above, the all-"*" placeholder row (originally
# row 28, designator "R10, R11") is now row 25.
$ws.Rows("25:25").Delete()

# Restore the user's last on-screen selection.
$ws.Range("E15").Select()
